# Update column F ("dSF") values for specific rows per repull of source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    11 = -1
    19 = -5
    22 = 0
    24 = -5
    25 = 5
    26 = -2
    27 = 1
    28 = 5
    29 = -8
    33 = 6
    37 = -4
    39 = 4
    43 = 4
    45 = 6
    46 = -2
    47 = -2
    55 = -1
    56 = -2
    58 = 0
    65 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
